$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Target cluster FAPs -> ECs (new cluster), refreshed metrics ---
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.520808
$ws.Range("H2").Value = 7.562424
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 22.41709
$ws.Range("N2").Value = 67.25127000000001
$ws.Range("O2").Value = 0.3988455747018376
$ws.Range("P2").Value = 0.3988455747018376
$ws.Range("Q2").Value = 56.50917980872001
$ws.Range("R2").Value = 508.58261827848
$ws.Range("S2").Value = 0.3988455747018376
$ws.Range("T2").Value = 0.3988455747018376

# --- Row 3: Target cluster sCs -> FAPs, refreshed metrics ---
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.520808
$ws.Range("H3").Value = 7.562424
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 16.78189033333333
$ws.Range("N3").Value = 50.345671
$ws.Range("O3").Value = 0.2985839238983091
$ws.Range("P3").Value = 0.2985839238983091
$ws.Range("Q3").Value = 42.30392340738933
$ws.Range("R3").Value = 380.7353106665039
$ws.Range("S3").Value = 0.2985839238983091
$ws.Range("T3").Value = 0.2985839238983091

# --- Row 4: new row, Target cluster sCs ---
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt2"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.520808
$ws.Range("H4").Value = 7.562424
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 17.00595566666667
$ws.Range("N4").Value = 51.017867
$ws.Range("O4").Value = 0.3025705013998533
$ws.Range("P4").Value = 0.3025705013998533
$ws.Range("Q4").Value = 42.86874909217867
$ws.Range("R4").Value = 385.818741829608
$ws.Range("S4").Value = 0.3025705013998533
$ws.Range("T4").Value = 0.3025705013998533
